$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B6").Value = 0.700000837771222
$ws.Range("C6").Value = 0.705029719932936
$ws.Range("D6").Value = 0.80010373602437
$ws.Range("E6").Value = 0.800087878992781
$ws.Range("F6").Value = 0.894985176165355
$ws.Range("G6").Value = 0.899996816599742
$ws.Range("C7").Value = 0.70518632778083
$ws.Range("D7").Value = 0.799842434518364
$ws.Range("E7").Value = 0.799387567606755
$ws.Range("F7").Value = 0.895146227880614
$ws.Range("G7").Value = 0.899997866572812
$ws.Range("B8").Value = 0.700002807471901
$ws.Range("C8").Value = 0.705009833694203
$ws.Range("D8").Value = 0.799936282470425
$ws.Range("E8").Value = 0.800015754951164
$ws.Range("F8").Value = 0.895025935817976
$ws.Range("B9").Value = 0.700003120908514
$ws.Range("C9").Value = 0.70510161321261
$ws.Range("D9").Value = 0.800186654178225
$ws.Range("E9").Value = 0.800281350966543
$ws.Range("F9").Value = 0.894973361308221
$ws.Range("G9").Value = 0.899999094428495
$ws.Range("C14").Value = 2.14349247503871
$ws.Range("D14").Value = 2.88586053978329
$ws.Range("E14").Value = 2.91716835666381
$ws.Range("F14").Value = 3.45922912691538
$ws.Range("G14").Value = 3.49993246406701
$ws.Range("C15").Value = 1.54055711441812
$ws.Range("D15").Value = 2.11544205582092
$ws.Range("E15").Value = 2.08594012344671
$ws.Range("F15").Value = 2.86242209502481
$ws.Range("G15").Value = 2.99991086842186
$ws.Range("C16").Value = 1.54084173711942
$ws.Range("D16").Value = 2.11745501779384
$ws.Range("E16").Value = 2.09003114090372
$ws.Range("F16").Value = 2.86086508112847
$ws.Range("G16").Value = 2.99997676691084
$ws.Range("C17").Value = 1.54108149620735
$ws.Range("D17").Value = 2.1124575663501
$ws.Range("E17").Value = 2.08305389764307
$ws.Range("F17").Value = 2.85124692755037
$ws.Range("G17").Value = 2.99999106774496
$ws.Range("C18").Value = 4034.27347442223
$ws.Range("D18").Value = 6686.7142658602
$ws.Range("E18").Value = 6684.70973116847
$ws.Range("F18").Value = 9317.09305417794
$ws.Range("C19").Value = 3266.72344181775
$ws.Range("D19").Value = 5385.31759111983
$ws.Range("E19").Value = 5387.46337253575
$ws.Range("F19").Value = 7485.3750789994
$ws.Range("C20").Value = 3108.01633241942
$ws.Range("D20").Value = 5105.51327062952
$ws.Range("E20").Value = 5104.40197156224
$ws.Range("F20").Value = 7112.61323808042
$ws.Range("C21").Value = 2915.68380919859
$ws.Range("D21").Value = 4791.62823733982
$ws.Range("E21").Value = 4790.03851118409
$ws.Range("F21").Value = 6673.10442852008
$ws.Range("C22").Value = 41.2561501054474
$ws.Range("D22").Value = 67.8774488630255
$ws.Range("E22").Value = 67.8540371297722
$ws.Range("F22").Value = 94.4568605872071
$ws.Range("B23").Value = 4.96925846828866
$ws.Range("C23").Value = 14.0972120615102
$ws.Range("D23").Value = 23.2033291619894
$ws.Range("E23").Value = 23.207353541631
$ws.Range("F23").Value = 32.2516335597981
$ws.Range("C24").Value = 6.90531462109457
$ws.Range("D24").Value = 11.4129504119885
$ws.Range("E24").Value = 11.4188741763382
$ws.Range("F24").Value = 15.881612246514
$ws.Range("G24").Value = 20.8674841376725
$ws.Range("B25").Value = 3.40328065734365
$ws.Range("C25").Value = 9.88730879459068
$ws.Range("D25").Value = 16.2549820575045
$ws.Range("E25").Value = 16.2633686483988
$ws.Range("F25").Value = 22.5937105650637
$ws.Range("B26").Value = 0.0130474105123126
$ws.Range("C26").Value = 0.0208506779534726
$ws.Range("D26").Value = 0.0347242895190912
$ws.Range("E26").Value = 0.0332063699709871
$ws.Range("F26").Value = 0.0577040210089048
$ws.Range("G26").Value = 0.158124990194769
$ws.Range("B27").Value = 0.800001325972683
$ws.Range("C27").Value = 0.814211623621882
$ws.Range("D27").Value = 0.888386454635636
$ws.Range("E27").Value = 0.891270378185676
$ws.Range("F27").Value = 0.945932121184867
$ws.Range("G27").Value = 0.949999459979907
$ws.Range("B28").Value = 26573.7665467139
$ws.Range("C28").Value = 43297.6913998287
$ws.Range("D28").Value = 68258.4681748415
$ws.Range("E28").Value = 67932.3969749865
$ws.Range("F28").Value = 95004.9119868568
$ws.Range("G28").Value = 122363.347750648
$ws.Range("B29").Value = 3.66512402306511
$ws.Range("C29").Value = 6.51928160484418
$ws.Range("D29").Value = 12.0462759431406
$ws.Range("E29").Value = 11.4749353844378
$ws.Range("F29").Value = 21.0600264554761
$ws.Range("G29").Value = 68.0549346040447
$ws.Range("B30").Value = 0.0129934819994118
$ws.Range("C30").Value = 0.0189588176502505
$ws.Range("D30").Value = 0.0328402598828854
$ws.Range("E30").Value = 0.0311366501978254
$ws.Range("F30").Value = 0.0567087304150561
$ws.Range("G30").Value = 0.195882116508292
$ws.Range("B31").Value = 0.800006142042167
$ws.Range("C31").Value = 0.814346333368092
$ws.Range("D31").Value = 0.888379952562697
$ws.Range("E31").Value = 0.891179698942495
$ws.Range("F31").Value = 0.945904907801852
$ws.Range("G31").Value = 0.949999251232065
$ws.Range("B32").Value = 16567.7719540081
$ws.Range("C32").Value = 27220.5502931098
$ws.Range("D32").Value = 44537.6921322833
$ws.Range("E32").Value = 43597.6650568061
$ws.Range("F32").Value = 66549.5069994244
$ws.Range("G32").Value = 89477.5633534774
$ws.Range("B33").Value = 2.75042826709312
$ws.Range("C33").Value = 5.0827313085502
$ws.Range("D33").Value = 9.75328325654722
$ws.Range("E33").Value = 9.20465801498353
$ws.Range("F33").Value = 17.6014844141599
$ws.Range("G33").Value = 54.5310533937057
$ws.Range("B34").Value = 0.0142792519117268
$ws.Range("C34").Value = 0.0221614525117038
$ws.Range("D34").Value = 0.0382341370875108
$ws.Range("E34").Value = 0.0363187175241403
$ws.Range("F34").Value = 0.0657783673656661
$ws.Range("G34").Value = 0.187668617401523
$ws.Range("B35").Value = 0.800009211058189
$ws.Range("C35").Value = 0.814148498574378
$ws.Range("D35").Value = 0.888337710853499
$ws.Range("E35").Value = 0.891162796071955
$ws.Range("F35").Value = 0.945855068794626
$ws.Range("G35").Value = 0.949999530738597
$ws.Range("B36").Value = 16206.7824859112
$ws.Range("C36").Value = 20980.2367193207
$ws.Range("D36").Value = 33173.3773029908
$ws.Range("E36").Value = 32469.1440004624
$ws.Range("F36").Value = 49139.1688712308
$ws.Range("G36").Value = 66714.4810995215
$ws.Range("B37").Value = 2.30566775632877
$ws.Range("C37").Value = 4.0868155582245
$ws.Range("D37").Value = 7.66524912492238
$ws.Range("E37").Value = 7.23804151442624
$ws.Range("F37").Value = 13.7425465072176
$ws.Range("G37").Value = 45.9506919811601
$ws.Range("B38").Value = 0.0160773458875501
$ws.Range("C38").Value = 0.024931790675803
$ws.Range("D38").Value = 0.0429471754905782
$ws.Range("E38").Value = 0.0407529955178683
$ws.Range("F38").Value = 0.0737081788563986
$ws.Range("G38").Value = 0.221902363953892
$ws.Range("B39").Value = 0.800000465331535
$ws.Range("C39").Value = 0.814350560180317
$ws.Range("D39").Value = 0.888610248974216
$ws.Range("E39").Value = 0.891444376008574
$ws.Range("F39").Value = 0.94602126312589
$ws.Range("G39").Value = 0.94999861201333
$ws.Range("B40").Value = 13966.2933007937
$ws.Range("C40").Value = 17836.6680338383
$ws.Range("D40").Value = 28008.9049724512
$ws.Range("E40").Value = 27382.4635954802
$ws.Range("F40").Value = 41389.8506809832
$ws.Range("G40").Value = 53864.2001116769
$ws.Range("B41").Value = 2.31046905241421
$ws.Range("C41").Value = 3.71146124221033
$ws.Range("D41").Value = 6.89203116789138
$ws.Range("E41").Value = 6.51134305209396
$ws.Range("F41").Value = 12.2987630808469
$ws.Range("G41").Value = 40.1646593037279
$ws.Range("B42").Value = 107.400316225466
$ws.Range("C42").Value = 367.249722445567
$ws.Range("D42").Value = 817.847074407423
$ws.Range("E42").Value = 768.415004254002
$ws.Range("F42").Value = 1555.31199161448
$ws.Range("G42").Value = 5547.53378840052
$ws.Range("B43").Value = 37.1657268328869
$ws.Range("C43").Value = 98.8821777998532
$ws.Range("D43").Value = 226.296445450132
$ws.Range("E43").Value = 211.080962133639
$ws.Range("F43").Value = 441.840262360131
$ws.Range("G43").Value = 1523.64965477712
$ws.Range("B44").Value = 9.43058458605051
$ws.Range("C44").Value = 38.9754910729192
$ws.Range("D44").Value = 87.4720225512977
$ws.Range("E44").Value = 81.6264786969934
$ws.Range("F44").Value = 170.138510664147
$ws.Range("G44").Value = 478.384233961743
$ws.Range("B45").Value = 16.7905191778637
$ws.Range("C45").Value = 50.1518996618518
$ws.Range("D45").Value = 111.998842878935
$ws.Range("E45").Value = 104.564557857329
$ws.Range("F45").Value = 216.485444670416
$ws.Range("G45").Value = 645.859953771358
$ws.Range("B46").Value = 86.6132569649369
$ws.Range("C46").Value = 289.4034148682
$ws.Range("D46").Value = 654.354812977798
$ws.Range("E46").Value = 613.315131915703
$ws.Range("F46").Value = 1257.14657552241
$ws.Range("G46").Value = 4565.06639238687
$ws.Range("B47").Value = 28.3593222415416
$ws.Range("C47").Value = 77.8236461455844
$ws.Range("D47").Value = 181.008823145205
$ws.Range("E47").Value = 168.498350875499
$ws.Range("F47").Value = 356.595796255891
$ws.Range("G47").Value = 1110.13780959192
$ws.Range("B48").Value = 8.46878948419991
$ws.Range("C48").Value = 30.6825809816439
$ws.Range("D48").Value = 69.968036862067
$ws.Range("E48").Value = 65.2376521942492
$ws.Range("F48").Value = 136.941496224943
$ws.Range("G48").Value = 362.291476565331
$ws.Range("B49").Value = 12.2159977022572
$ws.Range("C49").Value = 39.6121594815483
$ws.Range("D49").Value = 89.6195254042993
$ws.Range("E49").Value = 83.4271454342877
$ws.Range("F49").Value = 174.826431899219
$ws.Range("G49").Value = 514.894338854913
